$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text runs (preserve other run text, only change the targeted substring) ---

# A8: "Volume 32   Number  17" -> "Volume 32   Number  18"
$cellA8 = $ws.Range("A8")
$textA8 = $cellA8.Value2
$idx = $textA8.IndexOf("17")
$cellA8.Characters($idx + 1, 2).Text = "18"

# C9: "Report Covering the Week  4/21/2025  Through  4/27/2025"
#     -> "Report Covering the Week  4/28/2025  Through  5/4/2025"
$cellC9 = $ws.Range("C9")
$textC9 = $cellC9.Value2
$idx1 = $textC9.IndexOf("4/21/2025")
$cellC9.Characters($idx1 + 1, 9).Text = "4/28/2025"
$textC9b = $cellC9.Value2
$idx2 = $textC9b.IndexOf("4/27/2025")
$cellC9.Characters($idx2 + 1, 9).Text = "5/4/2025"

# --- Update weekly crime statistics table (rows 15-30) ---

# Row 15
$ws.Range("F15").Value2 = 2
$ws.Range("G15").Value2 = 2
$ws.Range("I15").Value2 = 12
$ws.Range("J15").Value2 = 16
$ws.Range("K15").Value2 = -25
$ws.Range("L15").Value2 = -25
$ws.Range("M15").Value2 = 140
$ws.Range("N15").Value2 = -42.857142857142

# Row 16
$ws.Range("C16").Value2 = 9
$ws.Range("D16").Value2 = 11
$ws.Range("E16").Value2 = -18.181818181818
$ws.Range("F16").Value2 = 38
$ws.Range("G16").Value2 = 36
$ws.Range("H16").Value2 = 5.555555555555
$ws.Range("I16").Value2 = 145
$ws.Range("J16").Value2 = 139
$ws.Range("K16").Value2 = 4.316546762589
$ws.Range("L16").Value2 = 12.403100775193
$ws.Range("M16").Value2 = 79.012345679012
$ws.Range("N16").Value2 = -63.104325699745

# Row 17
$ws.Range("C17").Value2 = 24
$ws.Range("D17").Value2 = 17
$ws.Range("E17").Value2 = 41.176470588235
$ws.Range("F17").Value2 = 59
$ws.Range("G17").Value2 = 66
$ws.Range("H17").Value2 = -10.60606060606
$ws.Range("I17").Value2 = 224
$ws.Range("J17").Value2 = 241
$ws.Range("K17").Value2 = -7.053941908713
$ws.Range("L17").Value2 = 13.705583756345
$ws.Range("M17").Value2 = 163.529411764706
$ws.Range("N17").Value2 = -24.324324324324

# Row 18
$ws.Range("C18").Value2 = 7
$ws.Range("D18").Value2 = 4
$ws.Range("E18").Value2 = 75
$ws.Range("G18").Value2 = 21
$ws.Range("H18").Value2 = 14.285714285714
$ws.Range("I18").Value2 = 88
$ws.Range("J18").Value2 = 97
$ws.Range("K18").Value2 = -9.278350515463
$ws.Range("L18").Value2 = 25.714285714285
$ws.Range("M18").Value2 = 114.634146341463
$ws.Range("N18").Value2 = -75.211267605633

# Row 19
$ws.Range("C19").Value2 = 9
$ws.Range("D19").Value2 = 11
$ws.Range("E19").Value2 = -18.181818181818
$ws.Range("F19").Value2 = 39
$ws.Range("G19").Value2 = 41
$ws.Range("H19").Value2 = -4.878048780487
$ws.Range("I19").Value2 = 171
$ws.Range("J19").Value2 = 200
$ws.Range("K19").Value2 = -14.5
$ws.Range("L19").Value2 = 12.5
$ws.Range("M19").Value2 = 137.5
$ws.Range("N19").Value2 = 66.019417475728

# Row 20
$ws.Range("D20").Value2 = 5
$ws.Range("E20").Value2 = 0
$ws.Range("F20").Value2 = 25
$ws.Range("G20").Value2 = 19
$ws.Range("H20").Value2 = 31.578947368421
$ws.Range("I20").Value2 = 89
$ws.Range("J20").Value2 = 81
$ws.Range("K20").Value2 = 9.876543209876
$ws.Range("L20").Value2 = -44.025157232704
$ws.Range("M20").Value2 = 117.073170731707
$ws.Range("N20").Value2 = -54.82233502538

# Row 21
$ws.Range("C21").Value2 = 55
$ws.Range("D21").Value2 = 49
$ws.Range("E21").Value2 = 12.244897959183
$ws.Range("F21").Value2 = 187
$ws.Range("G21").Value2 = 185
$ws.Range("H21").Value2 = 1.081081081081
$ws.Range("I21").Value2 = 729
$ws.Range("J21").Value2 = 775
$ws.Range("K21").Value2 = -5.935483870967
$ws.Range("L21").Value2 = 0.275103163686
$ws.Range("M21").Value2 = 120.909090909091
$ws.Range("N21").Value2 = -46.904588492352

# Row 23
$ws.Range("C23").Value2 = 6
$ws.Range("D23").Value2 = 3
$ws.Range("E23").Value2 = 100
$ws.Range("F23").Value2 = 24
$ws.Range("G23").Value2 = 25
$ws.Range("H23").Value2 = -4
$ws.Range("I23").Value2 = 92
$ws.Range("J23").Value2 = 132
$ws.Range("K23").Value2 = -30.30303030303
$ws.Range("L23").Value2 = -35.664335664335
$ws.Range("M23").Value2 = 48.387096774193

# Row 24
$ws.Range("C24").Value2 = 26
$ws.Range("D24").Value2 = 16
$ws.Range("E24").Value2 = 62.5
$ws.Range("G24").Value2 = 89
$ws.Range("H24").Value2 = 23.595505617977
$ws.Range("I24").Value2 = 410
$ws.Range("J24").Value2 = 380
$ws.Range("K24").Value2 = 7.894736842105
$ws.Range("L24").Value2 = 7.894736842105
$ws.Range("M24").Value2 = 65.991902834008

# Row 25
$ws.Range("C25").Value2 = 6
$ws.Range("D25").Value2 = 4
$ws.Range("E25").Value2 = 50
$ws.Range("F25").Value2 = 31
$ws.Range("G25").Value2 = 16
$ws.Range("H25").Value2 = 93.75
$ws.Range("I25").Value2 = 83
$ws.Range("J25").Value2 = 58
$ws.Range("K25").Value2 = 43.103448275862
$ws.Range("L25").Value2 = 5.06329113924

# Row 26
$ws.Range("C26").Value2 = 22
$ws.Range("D26").Value2 = 18
$ws.Range("E26").Value2 = 22.222222222222
$ws.Range("F26").Value2 = 78
$ws.Range("G26").Value2 = 73
$ws.Range("H26").Value2 = 6.849315068493
$ws.Range("I26").Value2 = 299
$ws.Range("J26").Value2 = 291
$ws.Range("K26").Value2 = 2.74914089347
$ws.Range("L26").Value2 = -22.337662337662
$ws.Range("M26").Value2 = 5.653710247349

# Row 27
$ws.Range("F27").Value2 = 2
$ws.Range("H27").Value2 = -33.333333333333
$ws.Range("I27").Value2 = 13
$ws.Range("J27").Value2 = 22
$ws.Range("K27").Value2 = -40.90909090909
$ws.Range("L27").Value2 = -43.478260869565

# Row 28
$ws.Range("C28").Value2 = 2
$ws.Range("D28").Value2 = 4
$ws.Range("G28").Value2 = 9
$ws.Range("H28").Value2 = -22.222222222222
$ws.Range("I28").Value2 = 19
$ws.Range("J28").Value2 = 27
$ws.Range("K28").Value2 = -29.629629629629
$ws.Range("L28").Value2 = -45.714285714285

# Row 29
$ws.Range("D29").Value2 = 1
$ws.Range("G29").Value2 = 4
$ws.Range("H29").Value2 = -75
$ws.Range("J29").Value2 = 7
$ws.Range("K29").Value2 = 0
$ws.Range("M29").Value2 = -53.333333333333
$ws.Range("N29").Value2 = -73.076923076923

# Row 30
$ws.Range("D30").Value2 = 1
$ws.Range("G30").Value2 = 4
$ws.Range("H30").Value2 = -75
$ws.Range("J30").Value2 = 7
$ws.Range("K30").Value2 = -14.285714285714
$ws.Range("M30").Value2 = -53.846153846153
$ws.Range("N30").Value2 = -76.923076923076
